# Applies the cryptocurrency price/volume refresh described in the commit diff.
# All target cells are plain text (t="inlineStr" in the original OOXML) even when
# their content looks numeric (e.g. "0.581", "34.687.36"), so each write forces a
# text number-format before assignment (otherwise Excel COM auto-coerces bare
# decimals into floating-point numbers, corrupting values like "34.99" into
# "34.990000000000002"), then clears the format again afterwards so the cell keeps
# its original (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", '34.687.36'),
    @("E2", '  +0.67%  '),
    @("D3", '1.821.40'),
    @("E3", '  +1.57%  '),
    @("E4", '  +0.10%  '),
    @("D5", '229.09'),
    @("E5", '  +1.13%  '),
    @("D6", '0.581'),
    @("E6", '  +4.44%  '),
    @("E7", '  +0.07%  '),
    @("D8", '34.99'),
    @("E8", '  +7.14%  '),
    @("E9", '  +1.72%  '),
    @("D10", '0.0699'),
    @("E10", '  +0.70%  '),
    @("D11", '0.0953'),
    @("E11", '  +0.37%  '),
    @("D12", '2.083.86'),
    @("E12", '  +1.51%  '),
    @("D13", '11.44'),
    @("E13", '  +3.53%  '),
    @("D14", '1.815.10'),
    @("E14", '  +0.14%  '),
    @("E15", '  +1.99%  '),
    @("D16", '34.659.56'),
    @("E16", '  +0.65%  '),
    @("D17", '4.37'),
    @("E17", '  +2.29%  '),
    @("D18", '69.58'),
    @("E18", '  +0.99%  '),
    @("B19", 'ShibaInu'),
    @("C19", 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'),
    @("D19", '0.0₃0802'),
    @("E19", '  +0.39%  '),
    @("B20", 'BitcoinCash'),
    @("C20", 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'),
    @("D20", '246.68'),
    @("E20", '  -0.17%  '),
    @("D21", '11.65'),
    @("E21", '  +4.44%  '),
    @("E22", '  +0.15%  '),
    @("D23", '4.21'),
    @("E23", '  +0.99%  '),
    @("D24", '173.07'),
    @("E24", '  +5.53%  '),
    @("E25", '  +1.55%  '),
    @("D26", '7.57'),
    @("E26", '  +4.51%  '),
    @("D27", '16.87'),
    @("E27", '  +2.21%  '),
    @("E28", '  +2.49%  '),
    @("E29", '  -0.10%  '),
    @("E30", '  +2.97%  '),
    @("E31", '  +1.88%  '),
    @("E32", '  +1.17%  '),
    @("E33", '  +1.29%  '),
    @("E34", '  +1.73%  '),
    @("D35", '1.402.76'),
    @("E35", '  -1.79%  '),
    @("B36", 'ImmutableX'),
    @("C36", 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'),
    @("D36", '0.682'),
    @("E36", '  +2.31%  '),
    @("B37", 'RenderToken'),
    @("C37", 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'),
    @("D37", '2.55'),
    @("E37", '  -1.05%  '),
    @("E38", '  +0.18%  '),
    @("E39", '  +0.29%  '),
    @("E40", '  +4.95%  '),
    @("D41", '83.42'),
    @("E41", '  -1.42%  '),
    @("D42", '0.955'),
    @("E42", '  +1.96%  '),
    @("E43", '  +0.15%  '),
    @("D44", '13.86'),
    @("E44", '  +2.80%  '),
    @("E45", '  +2.51%  '),
    @("E46", '  -1.99%  '),
    @("D47", '6.06'),
    @("E47", '  -0.67%  '),
    @("D48", '1.984.17'),
    @("E48", '  +1.76%  '),
    @("D49", '105.54'),
    @("E49", '  -0.01%  '),
    @("E50", '  -1.72%  '),
    @("E51", '  +0.05%  ')
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newVal = $pair[1]
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newVal
    $c.ClearFormats()
}
